$wb = $excel.ActiveWorkbook

# --- Category sheet: insert 5 new rows with new category names ---
$wsCat = $wb.Worksheets.Item("Category")
$wsCat.Rows("2:6").Insert()
$wsCat.Range("A2").Value = "Bìa cứng"
$wsCat.Range("A3").Value = "Bìa mềm"
$wsCat.Range("A4").Value = "Văn học"
$wsCat.Range("A5").Value = "Lịch sử"
$wsCat.Range("A6").Value = "Địa lý"

# --- Book sheet: add "categories" column G, and fill C24 ---
$wsBook = $wb.Worksheets.Item("Book")
$wsBook.Range("G1").Value = "categories"
$wsBook.Range("G2").Value = "1, 9"
$wsBook.Range("G3").Value = "2,11"
$wsBook.Range("G4").Value = "2,10"
$wsBook.Range("G5").Value = "2,3,12"
$wsBook.Range("G6").Value = "2,6"
$wsBook.Range("G7").Value = "2,10,14,4"
$wsBook.Range("G8").Value = "2,10"
$wsBook.Range("G9").Value = "1,3,10"
$wsBook.Range("G10").Value = "2,3"
$wsBook.Range("G11").Value = "2,7,9"
$wsBook.Range("G12").Value = "1, 6,9"
$wsBook.Range("G13").Value = "2, 10, 14"
$wsBook.Range("G14").Value = "2, 10, 11"
$wsBook.Range("G15").Value = "2,6,7"
$wsBook.Range("G16").Value = "1,9,14"
$wsBook.Range("G17").Value = "1,3"
$wsBook.Range("G18").Value = "1,12"
$wsBook.Range("G19").Value = "2,15,3"
$wsBook.Range("G20").Value = "2,5"
$wsBook.Range("G21").Value = "2,11"
$wsBook.Range("G22").Value = "2,3,10"
$wsBook.Range("G23").Value = "2,3"
$wsBook.Range("G24").Value = "2,8,4"
$wsBook.Range("G25").Value = "2,10,13"
$wsBook.Range("C24").Value = "Nhiều tác giả"

# --- Restore the UI selection state shown in the final workbook ---
[void]$wsCat.Range("C9").Select()
[void]$wsBook.Activate()
[void]$wsBook.Range("A2:XFD2").Select()

Write-Output "done"
